$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose displayed text must be preserved exactly as stored text
# (the sheet stores prices/percentages as text, not numbers),
# so force text format before assigning the values.
$textCells = @{
    'D2' = '337.77'
    'E2' = '2.70%'
    'D3' = '44.18'
    'E3' = '7.13%'
    'D4' = '5.775'
    'E4' = '2.61%'
    'D5' = '0.08340'
    'E5' = '1.81%'
    'D6' = '8.846'
    'E6' = '1.19%'
    'B7' = 'GateToken'
    'C7' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D7' = '4.526'
    'E7' = '0.93%'
    'B8' = 'FTXToken'
    'C8' = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
    'D8' = '1.973'
    'E8' = '-1.62%'
    'B9' = 'BTSEToken'
    'C9' = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    'D9' = '2.898'
    'E9' = '-2.97%'
    'B10' = 'MXToken'
    'C10' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D10' = '0.9456'
    'E10' = '2.72%'
    'B11' = 'LiechtensteinCryptoassetsExchange'
    'C11' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'D11' = '0.1246'
    'E11' = '-2.58%'
    'B12' = 'WazirX'
    'C12' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D12' = '0.1955'
    'E12' = '0.08%'
    'B13' = 'MandalaExchangeToken'
    'C13' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D13' = '0.09956'
    'E13' = '7.89%'
    'B14' = 'BitrueCoin'
    'C14' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D14' = '0.04525'
    'E14' = '16.43%'
    'B15' = 'BitMartToken'
    'C15' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D15' = '0.1068'
    'E15' = '1.02%'
    'B16' = 'BitForexToken'
    'C16' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D16' = '0.001302'
    'E16' = '0.42%'
    'B17' = 'TigerCash'
    'C17' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'D17' = '0.006082'
    'E17' = '-3.69%'
    'B18' = 'LEO'
    'C18' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D18' = '3.498'
    'E18' = '1.46%'
    'D20' = '8.761'
    'E20' = '6.50%'
    'D21' = '0.1371'
    'E21' = '-0.05%'
    'D22' = '0.2693'
    'E22' = '11.75%'
    'D23' = '0.04419'
    'E23' = '0.58%'
    'D24' = '0.001262'
    'E24' = '0.56%'
    'D25' = '0.004369'
    'E25' = '1.24%'
    'D26' = '0.0001262'
    'E26' = '5.05%'
    'D27' = '0.0003995'
    'D39' = '0.02802'
    'E39' = '0.50%'
    'D40' = '0.05819'
    'E40' = '7.64%'
    'D41' = '0.007950'
    'E41' = '3.29%'
    'D42' = '0.1432'
    'E42' = '1.05%'
    'D43' = '0.008968'
    'E43' = '0.16%'
    'D44' = '0.002123'
    'E44' = '-2.25%'
    'D45' = '0.008954'
    'E45' = '-22.02%'
    'D46' = '0.00007291'
    'E46' = '10.89%'
    'E47' = '0.11%'
    'D48' = '0.003196'
    'E48' = '-0.48%'
    'D49' = '0.002273'
    'E49' = '-0.29%'
    'D50' = '0.00002103'
    'E50' = '0.11%'
    'D51' = '0.0002003'
    'E51' = '0.11%'
}

foreach ($addr in $textCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textCells[$addr]
}
